# Applies the "tranverse gamma in pie dataset" edit to the MEDA results workbook.
# Adds a new gamma=0.4 column to the USOS_MNIST "relative performance" table
# (rows 59-67) and a new "gamma=0.1 transfer" column + diff column to the PIE
# transfer-accuracy table (rows 70-92), pushing the COIL summary table down by
# two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. USOS_MNIST "relative performance" table (rows 59-67): add gamma=0.4 col
# ------------------------------------------------------------------

# Row 59 header: move the row label from B59 into A59, put "MEDA" into B59
# (C59/D59 keep their existing gamma labels) and add the new "gamma=0.4"
# label in E59.
$ws.Range("A59").Value = "USOS_MNIST"
$ws.Range("B59").Value = "MEDA"
$ws.Range("E59").Value = "gamma=0.4"

# New gamma=0.4 raw accuracy figures for the two runs.
$ws.Range("E60").Value = 0.90222222222222204
$ws.Range("E61").Value = 0.73799999999999999

# Relative-performance formulas now span columns C:E (shared formula).
$ws.Range("C64").Formula = "=C60-`$B60"
$ws.Range("D64").Formula = "=D60-`$B60"
$ws.Range("E64").Formula = "=E60-`$B60"
$ws.Range("C65").Formula = "=C61-`$B61"
$ws.Range("D65").Formula = "=D61-`$B61"
$ws.Range("E65").Formula = "=E61-`$B61"

# Average relative performance row picks up the new columns too.
$ws.Range("D67").Formula = "=AVERAGE(D64:D65)"
$ws.Range("E67").Formula = "=AVERAGE(E64:E65)"

# ------------------------------------------------------------------
# 2. PIE transfer-accuracy table (rows 71-90): add a second ("gamma=0.1")
#    accuracy column plus its difference-from-baseline column.
# ------------------------------------------------------------------

$pieValues = @{
    71 = 0.39410681399631697
    72 = 0.45955882352941202
    73 = 0.63172123760889198
    74 = 0.33762254901960798
    75 = 0.44327731092437
    76 = 0.50919117647058798
    77 = 0.690297386602583
    78 = 0.36151960784313703
    79 = 0.45378151260504201
    80 = 0.46531614487415601
    81 = 0.71312706518474001
    82 = 0.43443627450980399
    83 = 0.69477791116446597
    84 = 0.71332105586249195
    85 = 0.77450980392156898
    86 = 0.52389705882352899
    87 = 0.39195678271308498
    88 = 0.35420503376304502
    89 = 0.44852941176470601
    90 = 0.49474316611595098
}

foreach ($row in 71..90) {
    $ws.Range("C$row").Value = $pieValues[$row]
    $ws.Range("D$row").Formula = "=C$row-B$row"
}

# ------------------------------------------------------------------
# 3. Make room for the new PIE average row: insert two rows before the old
#    row 92 (COIL header), so the COIL block shifts down to rows 94-101.
# ------------------------------------------------------------------

$ws.Rows("91:92").Insert()

# Clean up the stray formatted-but-empty cells the row insert drags along,
# then restore the lone formatted (percent, black font) blank cell at A91.
$ws.Range("A91:B92").Clear()
$ws.Range("A91").Font.Color = 0
$ws.Range("A91").NumberFormat = "0.0000%"

# New row 92: average of the PIE difference column.
$ws.Range("D92").Formula = "=AVERAGE(D71:D90)"

# ------------------------------------------------------------------
# 4. Restore the current selection (scrolled further down the sheet).
# ------------------------------------------------------------------

$ws.Range("J90").Select()
